# Forsythe Descendancy Numbered - apply "Updated and added Morrison family info" edits
# (updates provided by Gary Morrison)
#
# This script updates 11 cells in Sheet1, column A, that contain genealogy
# narrative text. The edits:
#   1. Expand the two Morrison/Coop-Livingston marriage detail lines (rows 110, 114)
#      with extra marriage/divorce dates.
#   2. Update the Charles Edward Forsythe III family block (rows 206-214) with
#      fuller names, birth dates, and birth locations.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A110").Value = "                     m. Rosa Livingston b. ? d. ? m. ? dv. Dec 1900 rm. May 1901 dv. about July 1901"
$ws.Range("A114").Value = "                     m. Rosa Mae Coop b. 29 Dec 1884 Davis, Iowa d. 7 Jul 1970 Paris, Davis, Iowa m. about Nov 1903"

$ws.Range("A206").Value = "                        9. Charles Edward Forsythe, III. b. 31 Jul 1955 Mt Pleasant, Westmoreland, Pennsylvania d. living in 2013"
$ws.Range("A207").Value = "                           m. Melinda Jane Steinman b. 5 Aug 1956 Canton, Stark, Ohio d. living in 2013 m. 11 Oct 1975"
$ws.Range("A208").Value = "                           10. Rochelle Forsythe b. 30 Jun 1978 Knoxville, Knox, Tennessee d. living in 2013"
$ws.Range("A209").Value = "                               m. Richard Thomas Mayes b. 7 Sep 1978 Cookeville, Putnam, Tennessee d. living in 2013 m. 21 May 2005"
$ws.Range("A210").Value = "                               11. Annaliese Madelyn Mayes b. 11 Jul 2008 Knoxville, Knox, Tennessee d. living in 2013"
$ws.Range("A211").Value = "                               11. Graham Benjamin Mayes b. 18 Jul 2011 Knoxville, Knox, Tennessee d. living in 2013"
$ws.Range("A212").Value = "                           10. Andrea Sue Forsythe b. 12 Aug 1982 Salem, Columbiana, Ohio d. living in 2013"
$ws.Range("A213").Value = "                           10. May Ling Forsythe b. 14 May 1989 Hong Kong, China d. living in 2013 (adopted)"
$ws.Range("A214").Value = "                           10. Jeremiah Jacob Forsythe b. 14 Jan 1992 Hong Kong, China d. living in 2013 (adopted)"
